$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: swap some labels around
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Data rows 2-7: 6x6 "one-hot" style grid with the single 1 moved to a new column each row
$data = @(
    @(0,1,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,0,1,0,0),
    @(0,0,0,0,1,0),
    @(0,0,1,0,0,0),
    @(1,0,0,0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
